# Swap the values of columns A, B, E, F, G, H, Q, R between rows 21 and 24.
# (The two observation records were re-ordered; columns C, D, I, J, K, N, P,
# S, T, U, V, W, Y, Z, AA, AB, AD, AE, AF, AG, AT, AW, AX, AY are identical
# between the two rows and therefore stay untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $columns) {
    $cell21 = $ws.Range($col + "21")
    $cell24 = $ws.Range($col + "24")

    $value21 = $cell21.Value2
    $value24 = $cell24.Value2

    $cell21.Value2 = $value24
    $cell24.Value2 = $value21
}
